# Daily attendance processing - 2025-11-14 13:52:38
# For every "Recorded By" cell (column G) that lists multiple comma-separated
# recorders and includes "System"/"system" among them, reverse the order of
# the comma-separated list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($null -eq $text -or $text -eq "") {
        continue
    }

    $hasComma = $text -like "*,*"
    $hasSystem = ($text -like "*System*") -or ($text -like "*system*")

    if ($hasComma -and $hasSystem) {
        $parts = $text -split ", "
        $reversed = $parts[-1..-($parts.Count)]
        $newText = $reversed -join ", "
        $ws.Cells.Item($r, 7).Value = $newText
    }
}
